# Auto commit - 10141735
# Adds 7 new rows (57-63) of maintenance-report data for 2025-10-14,
# updates the report title date, extends the print area / used range,
# and refreshes the sheet selection - mirroring the source XLSX diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Title cell (A1): bump "製表日期" (report-generated-on date)
# ---------------------------------------------------------------
$ws.Range("A1").Value = "萊爾富 工作統計表  篩選月份：202510   (  製表日期:2025-10-14  )"

# ---------------------------------------------------------------
# 2) Seed formatting for the new rows by cloning the zebra-striped
#    formats already used by rows 55 (shaded) and 56 (plain).
# ---------------------------------------------------------------
$ws.Range("A55:AK55").Copy()
$ws.Range("A57:AK57").PasteSpecial(-4122)
$ws.Range("A59:AK59").PasteSpecial(-4122)
$ws.Range("A61:AK61").PasteSpecial(-4122)
$ws.Range("A63:AK63").PasteSpecial(-4122)

$ws.Range("A56:AK56").Copy()
$ws.Range("A58:AK58").PasteSpecial(-4122)
$ws.Range("A60:AK60").PasteSpecial(-4122)
$ws.Range("A62:AK62").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 3) Row 57 - THILF03741 / 三重重新三
# ---------------------------------------------------------------
$ws.Cells.Item(57, 1).Value = 55
$ws.Cells.Item(57, 2).Value = '服務'
$ws.Cells.Item(57, 3).Value = 2025101747
$ws.Cells.Item(57, 6).Value = 3741
$ws.Cells.Item(57, 7).Value = '三重重新三'
$ws.Cells.Item(57, 8).Value = '新北市三重區'
$ws.Cells.Item(57, 17).Value = 'THILF03741'
$ws.Cells.Item(57, 18).Value = '新北一'
$ws.Cells.Item(57, 19).Value = '吳宗鴻'
$ws.Cells.Item(57, 20).Value = 1
$ws.Cells.Item(57, 21).Value = '已完工'
$ws.Cells.Item(57, 22).Value = '2025-10-14 11:52:06'
$ws.Cells.Item(57, 23).Value = '2025-10-14 09:30:00'
$ws.Cells.Item(57, 24).Value = '2025-10-14 10:00:00'
$ws.Cells.Item(57, 26).Value = 0.5
$ws.Cells.Item(57, 28).Value = '到場處理'
$ws.Cells.Item(57, 29).Value = 'PMQ4'
$ws.Cells.Item(57, 30).Value = 'O'
$ws.Cells.Item(57, 37).Value = 'O'

# ---------------------------------------------------------------
# 4) Row 58 - THILF02302 / 三重重安店
# ---------------------------------------------------------------
$ws.Cells.Item(58, 1).Value = 56
$ws.Cells.Item(58, 2).Value = '服務'
$ws.Cells.Item(58, 3).Value = 2025101750
$ws.Cells.Item(58, 6).Value = 2302
$ws.Cells.Item(58, 7).Value = '三重重安店'
$ws.Cells.Item(58, 8).Value = '新北市三重區'
$ws.Cells.Item(58, 17).Value = 'THILF02302'
$ws.Cells.Item(58, 18).Value = '新北一'
$ws.Cells.Item(58, 19).Value = '吳宗鴻'
$ws.Cells.Item(58, 20).Value = 1
$ws.Cells.Item(58, 21).Value = '已完工'
$ws.Cells.Item(58, 22).Value = '2025-10-14 11:52:40'
$ws.Cells.Item(58, 23).Value = '2025-10-14 10:10:00'
$ws.Cells.Item(58, 24).Value = '2025-10-14 10:30:00'
$ws.Cells.Item(58, 26).Value = 0.3
$ws.Cells.Item(58, 28).Value = '到場處理'
$ws.Cells.Item(58, 29).Value = 'PMQ4'
$ws.Cells.Item(58, 30).Value = 'O'
$ws.Cells.Item(58, 37).Value = 'O'

# ---------------------------------------------------------------
# 5) Row 59 - THILF04134 / 三重大同南
# ---------------------------------------------------------------
$ws.Cells.Item(59, 1).Value = 57
$ws.Cells.Item(59, 2).Value = '服務'
$ws.Cells.Item(59, 3).Value = 2025101751
$ws.Cells.Item(59, 6).Value = 4134
$ws.Cells.Item(59, 7).Value = '三重大同南'
$ws.Cells.Item(59, 8).Value = '新北市三重區'
$ws.Cells.Item(59, 17).Value = 'THILF04134'
$ws.Cells.Item(59, 18).Value = '新北一'
$ws.Cells.Item(59, 19).Value = '吳宗鴻'
$ws.Cells.Item(59, 20).Value = 1
$ws.Cells.Item(59, 21).Value = '已完工'
$ws.Cells.Item(59, 22).Value = '2025-10-14 11:53:15'
$ws.Cells.Item(59, 23).Value = '2025-10-14 11:00:00'
$ws.Cells.Item(59, 24).Value = '2025-10-14 11:30:00'
$ws.Cells.Item(59, 26).Value = 0.5
$ws.Cells.Item(59, 28).Value = '到場處理'
$ws.Cells.Item(59, 29).Value = 'PMQ4'
$ws.Cells.Item(59, 30).Value = 'O'
$ws.Cells.Item(59, 37).Value = 'O'

# ---------------------------------------------------------------
# 6) Row 60 - THILF02321 / 三重同安店
# ---------------------------------------------------------------
$ws.Cells.Item(60, 1).Value = 58
$ws.Cells.Item(60, 2).Value = '服務'
$ws.Cells.Item(60, 3).Value = 2025101760
$ws.Cells.Item(60, 6).Value = 2321
$ws.Cells.Item(60, 7).Value = '三重同安店'
$ws.Cells.Item(60, 8).Value = '新北市三重區'
$ws.Cells.Item(60, 17).Value = 'THILF02321'
$ws.Cells.Item(60, 18).Value = '新北一'
$ws.Cells.Item(60, 19).Value = '吳宗鴻'
$ws.Cells.Item(60, 20).Value = 1
$ws.Cells.Item(60, 21).Value = '已完工'
$ws.Cells.Item(60, 22).Value = '2025-10-14 12:16:39'
$ws.Cells.Item(60, 23).Value = '2025-10-14 11:50:00'
$ws.Cells.Item(60, 24).Value = '2025-10-14 12:16:00'
$ws.Cells.Item(60, 26).Value = 0.4
$ws.Cells.Item(60, 28).Value = '到場處理'
$ws.Cells.Item(60, 29).Value = 'PMQ4'
$ws.Cells.Item(60, 30).Value = 'O'
$ws.Cells.Item(60, 37).Value = 'O'

# ---------------------------------------------------------------
# 7) Row 61 - THILF04191 / 三重溪美店
# ---------------------------------------------------------------
$ws.Cells.Item(61, 1).Value = 59
$ws.Cells.Item(61, 2).Value = '服務'
$ws.Cells.Item(61, 3).Value = 2025101783
$ws.Cells.Item(61, 6).Value = 4191
$ws.Cells.Item(61, 7).Value = '三重溪美店'
$ws.Cells.Item(61, 8).Value = '新北市三重區'
$ws.Cells.Item(61, 17).Value = 'THILF04191'
$ws.Cells.Item(61, 18).Value = '新北一'
$ws.Cells.Item(61, 19).Value = '吳宗鴻'
$ws.Cells.Item(61, 20).Value = 1
$ws.Cells.Item(61, 21).Value = '已完工'
$ws.Cells.Item(61, 22).Value = '2025-10-14 15:07:53'
$ws.Cells.Item(61, 23).Value = '2025-10-14 14:39:00'
$ws.Cells.Item(61, 24).Value = '2025-10-14 15:00:00'
$ws.Cells.Item(61, 26).Value = 0.4
$ws.Cells.Item(61, 28).Value = '到場處理'
$ws.Cells.Item(61, 29).Value = 'PMQ4'
$ws.Cells.Item(61, 30).Value = 'O'
$ws.Cells.Item(61, 37).Value = 'O'

# ---------------------------------------------------------------
# 8) Row 62 - THILF04196 / 三重蝶愛店
# ---------------------------------------------------------------
$ws.Cells.Item(62, 1).Value = 60
$ws.Cells.Item(62, 2).Value = '服務'
$ws.Cells.Item(62, 3).Value = 2025101804
$ws.Cells.Item(62, 6).Value = 4196
$ws.Cells.Item(62, 7).Value = '三重蝶愛店'
$ws.Cells.Item(62, 8).Value = '新北市三重區'
$ws.Cells.Item(62, 17).Value = 'THILF04196'
$ws.Cells.Item(62, 18).Value = '新北一'
$ws.Cells.Item(62, 19).Value = '吳宗鴻'
$ws.Cells.Item(62, 20).Value = 1
$ws.Cells.Item(62, 21).Value = '已完工'
$ws.Cells.Item(62, 22).Value = '2025-10-14 15:51:44'
$ws.Cells.Item(62, 23).Value = '2025-10-14 15:30:00'
$ws.Cells.Item(62, 24).Value = '2025-10-14 15:51:00'
$ws.Cells.Item(62, 26).Value = 0.4
$ws.Cells.Item(62, 28).Value = '到場處理'
$ws.Cells.Item(62, 29).Value = 'PMQ4'
$ws.Cells.Item(62, 30).Value = 'O'
$ws.Cells.Item(62, 37).Value = 'O'

# ---------------------------------------------------------------
# 9) Row 63 - THILF0D194 / 北縣五華三店 (store #D194 as text)
# ---------------------------------------------------------------
$ws.Cells.Item(63, 1).Value = 61
$ws.Cells.Item(63, 2).Value = '服務'
$ws.Cells.Item(63, 3).Value = 2025101822
$ws.Cells.Item(63, 6).Value = 'D194'
$ws.Cells.Item(63, 7).Value = '北縣五華三店'
$ws.Cells.Item(63, 8).Value = '新北市三重區'
$ws.Cells.Item(63, 17).Value = 'THILF0D194'
$ws.Cells.Item(63, 18).Value = '新北一'
$ws.Cells.Item(63, 19).Value = '吳宗鴻'
$ws.Cells.Item(63, 20).Value = 1
$ws.Cells.Item(63, 21).Value = '已完工'
$ws.Cells.Item(63, 22).Value = '2025-10-14 16:22:05'
$ws.Cells.Item(63, 23).Value = '2025-10-14 16:00:00'
$ws.Cells.Item(63, 24).Value = '2025-10-14 16:21:00'
$ws.Cells.Item(63, 26).Value = 0.4
$ws.Cells.Item(63, 28).Value = '到場處理'
$ws.Cells.Item(63, 29).Value = 'PMQ4'
$ws.Cells.Item(63, 30).Value = 'O'
$ws.Cells.Item(63, 37).Value = 'O'

# ---------------------------------------------------------------
# 10) Word-wrap refresh: once a row has its "報修說明"(P) / "工作內容"
#     (AC) note filled in, those two cells pick up the wrapped
#     variant of the zebra style (P56/AC56 included - it gained the
#     wrap variant too even though its own text didn't change).
#     Row 63 (the very last / newest entry) has not been through
#     that refresh yet, so it intentionally keeps the non-wrapped style.
# ---------------------------------------------------------------
$wrapRows = @(56, 57, 58, 59, 60, 61, 62)
foreach ($r in $wrapRows) {
    $ws.Cells.Item($r, 16).WrapText = $true   # column P
    $ws.Cells.Item($r, 29).WrapText = $true   # column AC
}

# ---------------------------------------------------------------
# 11) Print area now spans through the newly added rows
# ---------------------------------------------------------------
$ws.PageSetup.PrintArea = '$A$1:$AK$63'

# ---------------------------------------------------------------
# 12) Selection moves to the new last row, matching the saved view
# ---------------------------------------------------------------
$ws.Range("A63").Select()
